$wb = $excel.ActiveWorkbook

# --- Metadata sheet updates ---
$meta = $wb.Worksheets.Item("Metadata")

# URL value
$meta.Range("B2").Value = "https://2rdoc.pt/ig/ios-lifestyle-medicine/ConceptMap/MindfulnessDiagnosticMap"

# Date value
$meta.Range("B8").Value = "2025-08-20T10:40:04+01:00"

# Source (ValueSet) value
$meta.Range("B15").Value = "https://2rdoc.pt/fhir/ValueSet/mindfulness-outcome-vs"

# --- Mapping Table 0 sheet updates ---
$map = $wb.Worksheets.Item("Mapping Table 0")

# Source (CodeSystem) value
$map.Range("A2").Value = "https://2rdoc.pt/fhir/CodeSystem/mindfulness-outcome-cs"
